$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.650468568325266
$ws.Cells.Item(2, 3).Value = 0.1255219744994989
$ws.Cells.Item(2, 5).Value = 0.2185565308620028
$ws.Cells.Item(2, 6).Value = 2.01452297117288
$ws.Cells.Item(2, 7).Value = 0.00247127532784681
$ws.Cells.Item(2, 9).Value = 0.9478811732533252
$ws.Cells.Item(2, 10).Value = 0.05414263894768467
$ws.Cells.Item(2, 11).Value = 0.3348554129341892
$ws.Cells.Item(2, 12).Value = 0.4808417783785188
$ws.Cells.Item(2, 14).Value = 1.813637775454717
$ws.Cells.Item(2, 15).Value = 3.569950141224936

$ws.Cells.Item(3, 2).Value = 0.6111677873460053
$ws.Cells.Item(3, 3).Value = 0.1250730213821498
$ws.Cells.Item(3, 5).Value = 0.217232952757989
$ws.Cells.Item(3, 6).Value = 2.013420968318755
$ws.Cells.Item(3, 7).Value = 0.002473448128908371
$ws.Cells.Item(3, 9).Value = 0.9548465628343976
$ws.Cells.Item(3, 10).Value = 0.05291071455588536
$ws.Cells.Item(3, 11).Value = 0.3018535448156001
$ws.Cells.Item(3, 12).Value = 0.4710515204788948
$ws.Cells.Item(3, 14).Value = 1.830701494907165
$ws.Cells.Item(3, 15).Value = 3.592369311787621

$ws.Cells.Item(4, 2).Value = 0.5872253557427314
$ws.Cells.Item(4, 3).Value = 0.1247973369156981
$ws.Cells.Item(4, 5).Value = 0.216501691645135
$ws.Cells.Item(4, 6).Value = 2.013686596389078
$ws.Cells.Item(4, 7).Value = 0.002474854122267539
$ws.Cells.Item(4, 9).Value = 0.9595254026845517
$ws.Cells.Item(4, 10).Value = 0.05214740876794366
$ws.Cells.Item(4, 11).Value = 0.2816247590335053
$ws.Cells.Item(4, 12).Value = 0.46523669614227
$ws.Cells.Item(4, 14).Value = 1.841724271868332
$ws.Cells.Item(4, 15).Value = 3.607784496648605

$ws.Cells.Item(5, 2).Value = 0.5775167512989867
$ws.Cells.Item(5, 3).Value = 0.1246849962944836
$ws.Cells.Item(5, 5).Value = 0.2162242399714813
$ws.Cells.Item(5, 6).Value = 2.014032112367232
$ws.Cells.Item(5, 7).Value = 0.002475445204254105
$ws.Cells.Item(5, 9).Value = 0.9615332325325561
$ws.Cells.Item(5, 10).Value = 0.05183463364134511
$ws.Cells.Item(5, 11).Value = 0.2733905660521998
$ws.Cells.Item(5, 12).Value = 0.4629167362139839
$ws.Cells.Item(5, 14).Value = 1.846353383873804
$ws.Cells.Item(5, 15).Value = 3.614481360111611

$ws.Cells.Item(6, 2).Value = 0.5759075741103459
$ws.Cells.Item(6, 3).Value = 0.1246663426826089
$ws.Cells.Item(6, 5).Value = 0.2161794121211749
$ws.Cells.Item(6, 6).Value = 2.014103825310045
$ws.Cells.Item(6, 7).Value = 0.002475544449515498
$ws.Cells.Item(6, 9).Value = 0.9618727430939416
$ws.Cells.Item(6, 10).Value = 0.05178259395823304
$ws.Cells.Item(6, 11).Value = 0.2720238558545702
$ws.Cells.Item(6, 12).Value = 0.4625345129168039
$ws.Cells.Item(6, 14).Value = 1.84713033489777
$ws.Cells.Item(6, 15).Value = 3.615618443539461

$ws.Cells.Item(7, 2).Value = 0.5870942262278618
$ws.Cells.Item(7, 3).Value = 0.1247958218255079
$ws.Cells.Item(7, 5).Value = 0.216497866564314
$ws.Cells.Item(7, 6).Value = 2.01369029500065
$ws.Cells.Item(7, 7).Value = 0.002474862020280049
$ws.Cells.Item(7, 9).Value = 0.9595520712851808
$ws.Cells.Item(7, 10).Value = 0.05214319752741403
$ws.Cells.Item(7, 11).Value = 0.2815136718714513
$ws.Cells.Item(7, 12).Value = 0.4652052070985917
$ws.Cells.Item(7, 14).Value = 1.841786145942478
$ws.Cells.Item(7, 15).Value = 3.607873132084563

$ws.Cells.Item(8, 2).Value = 0.6368789947517541
$ws.Cells.Item(8, 3).Value = 0.1253671862425101
$ws.Cells.Item(8, 5).Value = 0.2180833028298608
$ws.Cells.Item(8, 6).Value = 2.013947627447394
$ws.Cells.Item(8, 7).Value = 0.002472009621604252
$ws.Cells.Item(8, 9).Value = 0.9501994290386442
$ws.Cells.Item(8, 10).Value = 0.05371931087861626
$ws.Cells.Item(8, 11).Value = 0.3234695156217811
$ws.Cells.Item(8, 12).Value = 0.4774254637832627
$ws.Cells.Item(8, 14).Value = 1.819408140957714
$ws.Cells.Item(8, 15).Value = 3.577337997375651

$ws.Cells.Item(9, 2).Value = 0.7359736962249031
$ws.Cells.Item(9, 3).Value = 0.1264870709153705
$ws.Cells.Item(9, 5).Value = 0.2218356817062563
$ws.Cells.Item(9, 6).Value = 2.021916471368669
$ws.Cells.Item(9, 7).Value = 0.002466984056750266
$ws.Cells.Item(9, 9).Value = 0.9350464761932251
$ws.Cells.Item(9, 10).Value = 0.05675487509728683
$ws.Cells.Item(9, 11).Value = 0.4059993672380244
$ws.Cells.Item(9, 12).Value = 0.5029393190207401
$ws.Cells.Item(9, 14).Value = 1.779851084915384
$ws.Cells.Item(9, 15).Value = 3.530540360661064

$ws.Cells.Item(10, 2).Value = 0.809642788364954
$ws.Cells.Item(10, 3).Value = 0.1273091130642499
$ws.Cells.Item(10, 5).Value = 0.2249816295904488
$ws.Cells.Item(10, 6).Value = 2.032308585001303
$ws.Cells.Item(10, 7).Value = 0.002463634675680951
$ws.Cells.Item(10, 9).Value = 0.9258534554900457
$ws.Cells.Item(10, 10).Value = 0.05895104473497881
$ws.Cells.Item(10, 11).Value = 0.4667694911910303
$ws.Cells.Item(10, 12).Value = 0.5226200075242673
$ws.Cells.Item(10, 14).Value = 1.753421084668874
$ws.Cells.Item(10, 15).Value = 3.504123992136329

$ws.Cells.Item(11, 2).Value = 0.8433384940473161
$ws.Cells.Item(11, 3).Value = 0.1276828434241253
$ws.Cells.Item(11, 5).Value = 0.2264966541303366
$ws.Cells.Item(11, 6).Value = 2.038018947827013
$ws.Cells.Item(11, 7).Value = 0.00246218469995168
$ws.Cells.Item(11, 9).Value = 0.9220918489796759
$ws.Cells.Item(11, 10).Value = 0.0599426641885259
$ws.Cells.Item(11, 11).Value = 0.4944409650541104
$ws.Cells.Item(11, 12).Value = 0.5317746079766295
$ws.Cells.Item(11, 14).Value = 1.741968024462864
$ws.Cells.Item(11, 15).Value = 3.493834467952155

$ws.Cells.Item(12, 2).Value = 0.8561238646164497
$ws.Cells.Item(12, 3).Value = 0.1278243254187004
$ws.Cells.Item(12, 5).Value = 0.2270823582774071
$ws.Cells.Item(12, 6).Value = 2.040322380871729
$ws.Cells.Item(12, 7).Value = 0.002461646172599945
$ws.Cells.Item(12, 9).Value = 0.920727818430457
$ws.Cells.Item(12, 10).Value = 0.06031708483444831
$ws.Cells.Item(12, 11).Value = 0.5049228313090737
$ws.Cells.Item(12, 12).Value = 0.5352700230551761
$ws.Cells.Item(12, 14).Value = 1.737712982754159
$ws.Cells.Item(12, 15).Value = 3.490186335517279

$ws.Cells.Item(13, 2).Value = 0.8533691834740296
$ws.Cells.Item(13, 3).Value = 0.1277938567449297
$ws.Cells.Item(13, 5).Value = 0.2269556839444249
$ws.Cells.Item(13, 6).Value = 2.039820027278935
$ws.Cells.Item(13, 7).Value = 0.002461761685723791
$ws.Cells.Item(13, 9).Value = 0.9210189007285621
$ws.Cells.Item(13, 10).Value = 0.06023649511880791
$ws.Cells.Item(13, 11).Value = 0.5026652372375224
$ws.Cells.Item(13, 12).Value = 0.5345159473963434
$ws.Cells.Item(13, 14).Value = 1.738625735736269
$ws.Cells.Item(13, 15).Value = 3.49096098527275

$ws.Cells.Item(14, 2).Value = 0.844389847027827
$ws.Cells.Item(14, 3).Value = 0.1276944841226495
$ws.Cells.Item(14, 5).Value = 0.2265446003215317
$ws.Cells.Item(14, 6).Value = 2.03820562768486
$ws.Cells.Item(14, 7).Value = 0.002462140183944909
$ws.Cells.Item(14, 9).Value = 0.921978418905482
$ws.Cells.Item(14, 10).Value = 0.05997348981870942
$ws.Cells.Item(14, 11).Value = 0.4953032520924125
$ws.Cells.Item(14, 12).Value = 0.5320616027638891
$ws.Cells.Item(14, 14).Value = 1.741616316497856
$ws.Cells.Item(14, 15).Value = 3.493529358504816

$ws.Cells.Item(15, 2).Value = 0.8388930457974482
$ws.Cells.Item(15, 3).Value = 0.1276336098170816
$ws.Cells.Item(15, 5).Value = 0.2262943600494296
$ws.Cells.Item(15, 6).Value = 2.037235120220117
$ws.Cells.Item(15, 7).Value = 0.002462373396963511
$ws.Cells.Item(15, 9).Value = 0.9225740173350516
$ws.Cells.Item(15, 10).Value = 0.05981224982229705
$ws.Cells.Item(15, 11).Value = 0.4907942335127018
$ws.Cells.Item(15, 12).Value = 0.5305619847602969
$ws.Cells.Item(15, 14).Value = 1.743458810504306
$ws.Cells.Item(15, 15).Value = 3.495134892225508

$ws.Cells.Item(16, 2).Value = 0.8074443128443818
$ws.Cells.Item(16, 3).Value = 0.1272846836218022
$ws.Cells.Item(16, 5).Value = 0.2248843018233586
$ws.Cells.Item(16, 6).Value = 2.031955152737822
$ws.Cells.Item(16, 7).Value = 0.00246373091338762
$ws.Cells.Item(16, 9).Value = 0.9261077405758478
$ws.Cells.Item(16, 10).Value = 0.05888608922951732
$ws.Cells.Item(16, 11).Value = 0.4649615858156437
$ws.Cells.Item(16, 12).Value = 0.5220257717165282
$ws.Cells.Item(16, 14).Value = 1.754181041072012
$ws.Cells.Item(16, 15).Value = 3.504831187063957

$ws.Cells.Item(17, 2).Value = 0.7881978925869078
$ws.Cells.Item(17, 3).Value = 0.1270705651382542
$ws.Cells.Item(17, 5).Value = 0.2240407202154451
$ws.Cells.Item(17, 6).Value = 2.028967602718225
$ws.Cells.Item(17, 7).Value = 0.002464582541523319
$ws.Cells.Item(17, 9).Value = 0.9283831993897209
$ws.Cells.Item(17, 10).Value = 0.05831600571235285
$ws.Cells.Item(17, 11).Value = 0.4491205549289248
$ws.Cells.Item(17, 12).Value = 0.5168405815055337
$ws.Cells.Item(17, 14).Value = 1.760904797304631
$ws.Cells.Item(17, 15).Value = 3.511221901607854

$ws.Cells.Item(18, 2).Value = 0.7771451593285406
$ws.Cells.Item(18, 3).Value = 0.1269473896563014
$ws.Cells.Item(18, 5).Value = 0.2235634172460479
$ws.Cells.Item(18, 6).Value = 2.027341779931504
$ws.Cells.Item(18, 7).Value = 0.002465079313072131
$ws.Cells.Item(18, 9).Value = 0.929731550136033
$ws.Cells.Item(18, 10).Value = 0.05798741044019451
$ws.Cells.Item(18, 11).Value = 0.4400117738906886
$ws.Cells.Item(18, 12).Value = 0.5138772022522033
$ws.Cells.Item(18, 14).Value = 1.76482578204131
$ws.Cells.Item(18, 15).Value = 3.515060270722074

$ws.Cells.Item(19, 2).Value = 0.7734058896562033
$ws.Cells.Item(19, 3).Value = 0.1269056813872567
$ws.Cells.Item(19, 5).Value = 0.223403170114878
$ws.Cells.Item(19, 6).Value = 2.026807206791943
$ws.Cells.Item(19, 7).Value = 0.002465248704330542
$ws.Cells.Item(19, 9).Value = 0.9301948753330365
$ws.Cells.Item(19, 10).Value = 0.05787603427156895
$ws.Cells.Item(19, 11).Value = 0.436928154653458
$ws.Cells.Item(19, 12).Value = 0.5128771245631327
$ws.Cells.Item(19, 14).Value = 1.766162575748996
$ws.Cells.Item(19, 15).Value = 3.516387806565604

$ws.Cells.Item(20, 2).Value = 0.7902449224420138
$ws.Cells.Item(20, 3).Value = 0.1270933605774687
$ws.Cells.Item(20, 5).Value = 0.2241297035080017
$ws.Cells.Item(20, 6).Value = 2.029276057478597
$ws.Cells.Item(20, 7).Value = 0.002464491166760118
$ws.Cells.Item(20, 9).Value = 0.9281368781435368
$ws.Cells.Item(20, 10).Value = 0.05837676446261142
$ws.Cells.Item(20, 11).Value = 0.450806598860737
$ws.Cells.Item(20, 12).Value = 0.5173905886869221
$ws.Cells.Item(20, 14).Value = 1.760183488532164
$ws.Cells.Item(20, 15).Value = 3.510524771551502

$ws.Cells.Item(21, 2).Value = 0.8470266092058694
$ws.Cells.Item(21, 3).Value = 0.127723673477611
$ws.Cells.Item(21, 5).Value = 0.2266650205026863
$ws.Cells.Item(21, 6).Value = 2.03867599054297
$ws.Cells.Item(21, 7).Value = 0.002462028723834474
$ws.Cells.Item(21, 9).Value = 0.9216949459892874
$ws.Cells.Item(21, 10).Value = 0.06005077040027018
$ws.Cells.Item(21, 11).Value = 0.4974655611999026
$ws.Cells.Item(21, 12).Value = 0.5327817245216835
$ws.Cells.Item(21, 14).Value = 1.740735685355023
$ws.Cells.Item(21, 15).Value = 3.492768227487687

$ws.Cells.Item(22, 2).Value = 0.8842852011589457
$ws.Cells.Item(22, 3).Value = 0.1281353725093979
$ws.Cells.Item(22, 5).Value = 0.2283918932574807
$ws.Cells.Item(22, 6).Value = 2.045641352929835
$ws.Cells.Item(22, 7).Value = 0.00246048082854125
$ws.Cells.Item(22, 9).Value = 0.9178368592047192
$ws.Cells.Item(22, 10).Value = 0.06113850353263928
$ws.Cells.Item(22, 11).Value = 0.5279787374518889
$ws.Cells.Item(22, 12).Value = 0.5430082806116161
$ws.Cells.Item(22, 14).Value = 1.728503290168614
$ws.Cells.Item(22, 15).Value = 3.482610454276937

$ws.Cells.Item(23, 2).Value = 0.8643862731896945
$ws.Cells.Item(23, 3).Value = 0.127915666808164
$ws.Cells.Item(23, 5).Value = 0.2274638563635705
$ws.Cells.Item(23, 6).Value = 2.041848698621237
$ws.Cells.Item(23, 7).Value = 0.002461301361981005
$ws.Cells.Item(23, 9).Value = 0.9198637887482732
$ws.Cells.Item(23, 10).Value = 0.06055854402767125
$ws.Cells.Item(23, 11).Value = 0.5116917481567498
$ws.Cells.Item(23, 12).Value = 0.5375349248566579
$ws.Cells.Item(23, 14).Value = 1.734988221724409
$ws.Cells.Item(23, 15).Value = 3.487899471501805

$ws.Cells.Item(24, 2).Value = 0.7893194219000463
$ws.Cells.Item(24, 3).Value = 0.127083054995019
$ws.Cells.Item(24, 5).Value = 0.2240894502238504
$ws.Cells.Item(24, 6).Value = 2.029136319235931
$ws.Cells.Item(24, 7).Value = 0.002464532454944649
$ws.Cells.Item(24, 9).Value = 0.9282481148827664
$ws.Cells.Item(24, 10).Value = 0.05834929806743006
$ws.Cells.Item(24, 11).Value = 0.4500443432185648
$ws.Cells.Item(24, 12).Value = 0.5171418754277113
$ws.Cells.Item(24, 14).Value = 1.760509419635669
$ws.Cells.Item(24, 15).Value = 3.510839432373814

$ws.Cells.Item(25, 2).Value = 0.7090120173540981
$ws.Cells.Item(25, 3).Value = 0.126184213574021
$ws.Cells.Item(25, 5).Value = 0.2207520121553372
$ws.Cells.Item(25, 6).Value = 2.018962931280242
$ws.Cells.Item(25, 7).Value = 0.002468283147464557
$ws.Cells.Item(25, 9).Value = 0.9388048512093086
$ws.Cells.Item(25, 10).Value = 0.05593962528680407
$ws.Cells.Item(25, 11).Value = 0.3836476438274872
$ws.Cells.Item(25, 12).Value = 0.4958721639302155
$ws.Cells.Item(25, 14).Value = 1.790089676678928
$ws.Cells.Item(25, 15).Value = 3.541800821025717
